# Resultados SETAR, ARMA y ARIMA
# Refresh the forecasting-method results table (columns F:N, rows 1-25):
# the method columns are re-ordered (alphabetically) and re-populated with
# the updated simulation output values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a 9-column (F:N) x 25-row (1:25) array: row 1 = new headers,
# rows 2-25 = new simulation results for each forecasting method.
$arr = New-Object "object[,]" 25,9

$arr[0,0] = "Block Bootstrapping"
$arr[0,1] = "Sieve Bootstrap"
$arr[0,2] = "LSPM"
$arr[0,3] = "LSPMW"
$arr[0,4] = "AREPD"
$arr[0,5] = "MCPS"
$arr[0,6] = "AV-MCPS"
$arr[0,7] = "DeepAR"
$arr[0,8] = "EnCQR-LSTM"

$arr[1,0] = 1.226924623394829
$arr[1,1] = 0.2244562714330278
$arr[1,2] = 0.8158269524473714
$arr[1,3] = 0.7474871919534316
$arr[1,4] = 1.052634119590127
$arr[1,5] = 0.3741561748674948
$arr[1,6] = 0.4811661845849273
$arr[1,7] = 0.1793823735526282
$arr[1,8] = 1.19174488844377

$arr[2,0] = 1.778796605821918
$arr[2,1] = 0.7150802718271987
$arr[2,2] = 0.6252332625047938
$arr[2,3] = 1.547999074363817
$arr[2,4] = 1.658752776116569
$arr[2,5] = 0.5093420019891721
$arr[2,6] = 0.6665048839020538
$arr[2,7] = 0.5507379707770272
$arr[2,8] = 1.295481679568396

$arr[3,0] = 1.386062675709015
$arr[3,1] = 0.2332915509516632
$arr[3,2] = 0.3814740306092502
$arr[3,3] = 0.975909184646331
$arr[3,4] = 1.204785416641957
$arr[3,5] = 0.927660621597965
$arr[3,6] = 0.7854608458592057
$arr[3,7] = 0.2541262744218593
$arr[3,8] = 1.240040281390922

$arr[4,0] = 0.5349509838935558
$arr[4,1] = 1.463878511254781
$arr[4,2] = 1.847545942000418
$arr[4,3] = 0.3470519864326069
$arr[4,4] = 0.4025449673880657
$arr[4,5] = 1.517835474385723
$arr[4,6] = 1.01641162707222
$arr[4,7] = 1.624204504953172
$arr[4,8] = 1.058161099876447

$arr[5,0] = 0.4872451287413546
$arr[5,1] = 0.2131834752240013
$arr[5,2] = 1.607550714455432
$arr[5,3] = 0.5222527728341866
$arr[5,4] = 0.4554449993668057
$arr[5,5] = 0.1948922942347683
$arr[5,6] = 0.3432584831700483
$arr[5,7] = 0.6742657760268511
$arr[5,8] = 1.046586658523247

$arr[6,0] = 0.5601735230247373
$arr[6,1] = 0.3935683052506395
$arr[6,2] = 0.2235438256374703
$arr[6,3] = 0.3246357140290516
$arr[6,4] = 0.3986377338024754
$arr[6,5] = 0.3748093418793953
$arr[6,6] = 0.5137295931781937
$arr[6,7] = 0.184225657456271
$arr[6,8] = 1.099638667864436

$arr[7,0] = 0.5611035176912833
$arr[7,1] = 0.3773642457269758
$arr[7,2] = 0.2428709598684402
$arr[7,3] = 0.8871943794329584
$arr[7,4] = 0.5963136792778903
$arr[7,5] = 0.1512838697631637
$arr[7,6] = 0.2980657415593938
$arr[7,7] = 0.5898021793249799
$arr[7,8] = 1.039773028100319

$arr[8,0] = 0.8106682412194928
$arr[8,1] = 0.4661031590940097
$arr[8,2] = 0.8530713121210909
$arr[8,3] = 1.55580164840771
$arr[8,4] = 0.9600656753207362
$arr[8,5] = 1.112222296724392
$arr[8,6] = 0.3058378795468658
$arr[8,7] = 0.653731136370021
$arr[8,8] = 1.002511163976704

$arr[9,0] = 0.4988686177045323
$arr[9,1] = 1.155972343456121
$arr[9,2] = 0.53835975309544
$arr[9,3] = 0.3054037869163848
$arr[9,4] = 0.3914477096344753
$arr[9,5] = 1.150177668785864
$arr[9,6] = 1.149310659448664
$arr[9,7] = 0.7595436930053753
$arr[9,8] = 1.154715860363411

$arr[10,0] = 1.273663400352505
$arr[10,1] = 1.596601270325984
$arr[10,2] = 0.4890885719048976
$arr[10,3] = 2.361199590800857
$arr[10,4] = 1.514595684089588
$arr[10,5] = 1.382924932715614
$arr[10,6] = 1.18127627652726
$arr[10,7] = 1.343207018778314
$arr[10,8] = 0.9772355457876221

$arr[11,0] = 2.573203820909429
$arr[11,1] = 0.9176939952677238
$arr[11,2] = 3.182120400824687
$arr[11,3] = 3.863597783485248
$arr[11,4] = 2.882308576305129
$arr[11,5] = 1.64247811206229
$arr[11,6] = 1.519823503770773
$arr[11,7] = 1.575914032818821
$arr[11,8] = 1.474763295336273

$arr[12,0] = 3.040571943127722
$arr[12,1] = 0.9391122800233652
$arr[12,2] = 2.030584084822666
$arr[12,3] = 4.618486632314057
$arr[12,4] = 3.618470781169586
$arr[12,5] = 1.389423689285336
$arr[12,6] = 1.277891755814129
$arr[12,7] = 1.225247700571236
$arr[12,8] = 1.906174753729765

$arr[13,0] = 0.6487996839784451
$arr[13,1] = 0.4633062315487987
$arr[13,2] = 1.204519599740304
$arr[13,3] = 1.221712425855758
$arr[13,4] = 0.5421522313135088
$arr[13,5] = 1.033618436045697
$arr[13,6] = 1.163085633753081
$arr[13,7] = 0.2930731269332834
$arr[13,8] = 0.711961918205156

$arr[14,0] = 0.4940918778199035
$arr[14,1] = 0.2448425912944219
$arr[14,2] = 0.3380705489779425
$arr[14,3] = 0.9559678984815889
$arr[14,4] = 0.3540881160232546
$arr[14,5] = 0.3649411508335579
$arr[14,6] = 0.6628382763493815
$arr[14,7] = 0.3099987478288683
$arr[14,8] = 0.7379913101857172

$arr[15,0] = 0.4658038999891121
$arr[15,1] = 0.5487206595077961
$arr[15,2] = 0.8244278837418509
$arr[15,3] = 0.2710644823162839
$arr[15,4] = 0.6078709855481128
$arr[15,5] = 0.9617461120447737
$arr[15,6] = 0.392061953583276
$arr[15,7] = 0.5905378715551792
$arr[15,8] = 0.9152039536424872

$arr[16,0] = 0.4743502681473842
$arr[16,1] = 0.2434238848972567
$arr[16,2] = 0.4930468661545232
$arr[16,3] = 0.2673069223247122
$arr[16,4] = 0.5785197825753251
$arr[16,5] = 0.6495970299591682
$arr[16,6] = 0.6146198522687679
$arr[16,7] = 0.226154050337088
$arr[16,8] = 0.911205068968987

$arr[17,0] = 0.3028028632463684
$arr[17,1] = 0.2205408325185787
$arr[17,2] = 0.4992768748766976
$arr[17,3] = 0.5108078384372974
$arr[17,4] = 0.1834398199932399
$arr[17,5] = 0.1602494064900193
$arr[17,6] = 0.4445817169114058
$arr[17,7] = 0.2469306161443889
$arr[17,8] = 0.7987872126325608

$arr[18,0] = 0.5890968561237072
$arr[18,1] = 0.4477649677359817
$arr[18,2] = 0.590715656470006
$arr[18,3] = 1.127253449116638
$arr[18,4] = 0.4684860701125651
$arr[18,5] = 0.7832444312157071
$arr[18,6] = 0.7287586519836391
$arr[18,7] = 0.5547888759092583
$arr[18,8] = 0.7230735777736721

$arr[19,0] = 0.3373095206161351
$arr[19,1] = 0.2520063687550229
$arr[19,2] = 0.2401938483972749
$arr[19,3] = 0.5835075270528773
$arr[19,4] = 0.1952746629641808
$arr[19,5] = 0.2164191063252457
$arr[19,6] = 0.5297971593036553
$arr[19,7] = 0.268514464372243
$arr[19,8] = 0.7861212864965388

$arr[20,0] = 0.2876922893353361
$arr[20,1] = 0.2469432477824522
$arr[20,2] = 0.3435876395138338
$arr[20,3] = 0.3701533589031113
$arr[20,4] = 0.2082085681113187
$arr[20,5] = 0.2707113529280156
$arr[20,6] = 0.219769629620479
$arr[20,7] = 0.2616272452617234
$arr[20,8] = 0.8318061264422902

$arr[21,0] = 0.3164848642364121
$arr[21,1] = 0.3315118335905223
$arr[21,2] = 0.2551805011200191
$arr[21,3] = 0.2908682109888899
$arr[21,4] = 0.2998718948734481
$arr[21,5] = 0.3508442256111938
$arr[21,6] = 0.2173640992241089
$arr[21,7] = 0.3082299777084325
$arr[21,8] = 0.8621505881227747

$arr[22,0] = 0.3001852949098409
$arr[22,1] = 0.2648657154541696
$arr[22,2] = 0.3067764406955029
$arr[22,3] = 0.2783892953914174
$arr[22,4] = 0.3248967465151147
$arr[22,5] = 0.3528663314700787
$arr[22,6] = 0.405414146842329
$arr[22,7] = 0.2397212679579679
$arr[22,8] = 0.8679691708360295

$arr[23,0] = 0.8354785171735373
$arr[23,1] = 0.9337172020316944
$arr[23,2] = 0.9913100160409892
$arr[23,3] = 1.489645264773894
$arr[23,4] = 0.7679697640885021
$arr[23,5] = 1.405838273486609
$arr[23,6] = 0.9829086219696229
$arr[23,7] = 0.8646218897037269
$arr[23,8] = 0.7052299570010334

$arr[24,0] = 0.314631890779088
$arr[24,1] = 0.3671918221886344
$arr[24,2] = 0.2618095568138273
$arr[24,3] = 0.5248030286970484
$arr[24,4] = 0.1797737190646329
$arr[24,5] = 0.4037993018297222
$arr[24,6] = 0.3877028253937319
$arr[24,7] = 0.4255724429691736
$arr[24,8] = 0.796025473385392

$ws.Range("F1:N25").Value = $arr
